$wb = $excel.ActiveWorkbook

# Overview sheet: the dc7c6889...md row (row 3) status moves from
# "Handed back: in sync with en-US" to "Ready for handoff" for both
# the zh-cn and de-de columns.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: same status update for the dc7c6889 row (row 3), plus
# the Latest Handoff Datetime for the 1942fc03 batch (rows 2 and 3,
# which share the same handoff timestamp) advances from 11:38:57 to
# 11:41:15.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D2").Value = "2016-01-28 11:41:15"
$wsZhCn.Range("D3").Value = "2016-01-28 11:41:15"

# de-de sheet: same status update for the dc7c6889 row (row 3), plus
# the Latest Handoff Datetime advances from 11:39:10 to 11:41:28.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D2").Value = "2016-01-28 11:41:28"
$wsDeDe.Range("D3").Value = "2016-01-28 11:41:28"
